$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '29.409.49'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '1.912.68'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = '  +0.79%  '

$ws.Range("D5").Value = "'325.07"
$ws.Range("E5").Value = '  +0.59%  '

$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("D7").Value = "'0.4824"
$ws.Range("E7").Value = '  +1.51%  '

$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = '  +0.35%  '

$ws.Range("D9").Value = "'0.08203"
$ws.Range("E9").Value = '  +2.17%  '

$ws.Range("E10").Value = '  +1.87%  '

$ws.Range("D11").Value = "'23.49"
$ws.Range("E11").Value = '  +0.54%  '

$ws.Range("D12").Value = '1.932.45'
$ws.Range("E12").Value = '  +2.83%  '

$ws.Range("D13").Value = "'6.050"
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").Value = "'7.210"
$ws.Range("E14").Value = '  +2.20%  '

$ws.Range("D15").Value = "'91.04"
$ws.Range("E15").Value = '  +1.98%  '

$ws.Range("D16").Value = "'0.06808"
$ws.Range("E16").Value = '  +2.24%  '

$ws.Range("E17").Value = '  +0.82%  '

$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").Value = "'17.70"
$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").Value = '29.412.02'
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = "'5.621"
$ws.Range("E22").Value = '  +1.73%  '

$ws.Range("E23").Value = '  +0.87%  '

$ws.Range("D24").Value = "'2.179"
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").Value = '2.150.82'
$ws.Range("E25").Value = '  -1.90%  '

$ws.Range("D26").Value = "'6.583"
$ws.Range("E26").Value = '  +10.59%  '

$ws.Range("D27").Value = "'155.78"
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("E28").Value = '  +1.38%  '

$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("D30").Value = "'120.30"
$ws.Range("E30").Value = '  +2.04%  '

$ws.Range("D31").Value = "'1.020"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").Value = "'0.09562"
$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("D33").Value = "'5.601"
$ws.Range("E33").Value = '  +4.70%  '

$ws.Range("D34").Value = "'3.550"
$ws.Range("E34").Value = '  +0.73%  '

$ws.Range("D35").Value = "'1.367"
$ws.Range("E35").Value = '  -1.03%  '

$ws.Range("D36").Value = "'0.02284"
$ws.Range("E36").Value = '  +1.62%  '

$ws.Range("D37").Value = "'0.06108"
$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("D38").Value = "'1.177"
$ws.Range("E38").Value = '  +0.38%  '

$ws.Range("D39").Value = "'8.062"
$ws.Range("E39").Value = '  +2.28%  '

$ws.Range("D40").Value = "'0.5969"
$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("D41").Value = "'10.82"
$ws.Range("E41").Value = '  +7.16%  '

$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'2.416"
$ws.Range("E43").Value = '  +1.73%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = "'1.279"
$ws.Range("E44").Value = '  +1.37%  '

$ws.Range("D45").Value = "'0.07621"
$ws.Range("E45").Value = '  -1.20%  '

$ws.Range("D46").Value = "'12.39"
$ws.Range("E46").Value = '  +1.14%  '

$ws.Range("D47").Value = "'0.5574"
$ws.Range("E47").Value = '  +1.37%  '

$ws.Range("D48").Value = "'1.954"
$ws.Range("E48").Value = '  +1.83%  '

$ws.Range("D49").Value = "'117.72"
$ws.Range("E49").Value = '  +3.93%  '

$ws.Range("D50").Value = "'2.421"
$ws.Range("E50").Value = '  +3.82%  '

$ws.Range("E51").Value = '  +0.88%  '
